$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; I = -68.395700000000005; J = -8.5576000000000008; K = -13.580299999999999 }
    @{ Row = 3; I = -67.773200000000003; J = -22.257400000000001; K = -7.6257000000000001 }
    @{ Row = 5; I = -66.544899999999998; J = 2.7523; K = -15.3057 }
    @{ Row = 6; I = -61.795699999999997; J = 15.499599999999999; K = -20.139800000000001 }
    @{ Row = 7; I = -67.751099999999994; J = -21.2439; K = -7.0552000000000001 }
    @{ Row = 8; I = -31.254000000000001; J = -11.762700000000001; K = -4.7976999999999999 }
    @{ Row = 9; I = -31.254000000000001; J = -11.762700000000001; K = -4.7976999999999999 }
    @{ Row = 10; I = -31.254000000000001; J = -11.762700000000001; K = -4.7976999999999999 }
    @{ Row = 11; I = -31.254000000000001; J = -11.762700000000001; K = -4.7976999999999999 }
    @{ Row = 12; I = -58.363599999999998; J = 16.240300000000001; K = -17.109100000000002 }
    @{ Row = 13; I = -66.228700000000003; J = -30.240100000000002; K = -3.3279999999999998 }
    @{ Row = 14; I = -67.751099999999994; J = -21.2439; K = -7.0552000000000001 }
    @{ Row = 15; I = -67.751099999999994; J = -21.2439; K = -7.0552000000000001 }
    @{ Row = 16; I = -67.860500000000002; J = -6.8452000000000002; K = -10.079499999999999 }
    @{ Row = 17; I = -58.363599999999998; J = 16.240300000000001; K = -17.109100000000002 }
    @{ Row = 18; I = -58.363599999999998; J = 16.240300000000001; K = -17.109100000000002 }
    @{ Row = 19; I = -31.9497; J = -12.861599999999999; K = -2.8336999999999999 }
    @{ Row = 20; I = -31.254000000000001; J = -11.762700000000001; K = -4.7976999999999999 }
    @{ Row = 21; I = -63.590200000000003; J = -2.6019999999999999; K = -7.5965999999999996 }
    @{ Row = 22; I = -63.590200000000003; J = -2.6019999999999999; K = -7.5965999999999996 }
    @{ Row = 23; I = -58.363599999999998; J = 16.240300000000001; K = -17.109100000000002 }
    @{ Row = 24; I = -58.363599999999998; J = 16.240300000000001; K = -17.109100000000002 }
    @{ Row = 25; I = -67.751099999999994; J = -21.2439; K = -7.0552000000000001 }
    @{ Row = 26; I = -31.254000000000001; J = -11.762700000000001; K = -4.7976999999999999 }
    @{ Row = 27; I = -58.363599999999998; J = 16.240300000000001; K = -17.109100000000002 }
    @{ Row = 28; I = -63.590200000000003; J = -2.6019999999999999; K = -7.5965999999999996 }
    @{ Row = 29; I = -31.254000000000001; J = -11.762700000000001; K = -4.7976999999999999 }
    @{ Row = 30; I = -31.9497; J = -12.861599999999999; K = -2.8336999999999999 }
    @{ Row = 31; I = -63.590200000000003; J = -2.6019999999999999; K = -7.5965999999999996 }
    @{ Row = 32; I = -31.9497; J = -12.861599999999999; K = -2.8336999999999999 }
    @{ Row = 33; I = -32.891800000000003; J = -14.9907; K = -1.5783 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
    $ws.Cells.Item($u.Row, 11).Value = $u.K
}

$ws.Range("L4").Select()
